$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bedingungszuweisung")

$rowsData = @"
Code	Bedingung
AB19E	1
AB31R	1
AC12E	1
AH10Y	1
AK12I	2
AM21H	2
AU09R	3
BC05R	1
BM15R	2
BR11A	1
BU07E	1
CC09K	2
CC10A	1
CC21W	1
CE10B	3
CJ10Z	2
CK11R	2
CT29K	1
CU24R	2
DC28D	2
DI21I	2
DK11R	1
DL16R	2
DN28H	2
DT04N	2
EA09H	1
EE08O	1
EG25R	1
EI14R	1
EM07N	3
EM16N	1
FA12L	2
FA23L	1
FB04R	2
FK22R	1
HE12B	2
HT08S	3
II08L	2
IM09E	1
JA24N	3
JA26L	3
JC05N	1
JE03S	1
JE15A	3
JG04G	2
JK12R	1
JK20G	2
JK25L	1
JM07ß	2
JM22C	2
JS26S	3
JU01F	1
KA17A	2
KA19E	2
KB16K	1
KD16N	1
KE17E	2
KL19T	1
KL24U	1
KM23K	1
KS11H	2
KS28N	2
KU08D	1
LA09L	2
LA23T	1
LA26R	2
LB21E	1
LC05H	1
LE04A	1
LG27Z	3
LI20Z	1
LJ16R	1
LK24L	3
LM18T	2
LP06T	2
LP25R	2
LS13L	1
LW28E	2
LW28R	2
MC04A	2
MD22R	1
MD26R	2
ME06K	1
ME18L	2
MG08K	1
MH14A	1
MH25K	2
MI15D	2
MS09E	1
MS23B	3
MS28B	2
NM11M	1
OB20R	2
PA13R	2
PG11T	1
PS06G	3
RA03A	2
RF17R	2
RM13N	1
SB11Z	1
SC10R	2
SD30K	2
SE06H	2
SE12O	2
SG29T	3
SK12K	2
SL10R	2
SR08R	3
SR19R	2
SS11H	1
SS21A	1
SS22K	2
SU29N	1
SV03R	3
TB29R	1
TM01F	2
TM30L	2
UH09L	1
VA30T	1
VE19A	2
VS09S	2
WM13K	2
"@

$lines = $rowsData -split "`n"
$rowNum = 0
foreach ($line in $lines) {
    $line = $line.Trim("`r")
    if ($line -eq "") { continue }
    $rowNum++
    $parts = $line -split "`t"
    $code = $parts[0]
    $val = $parts[1]
    $ws.Cells.Item($rowNum, 1).Value = $code
    if ($rowNum -eq 1) {
        $ws.Cells.Item($rowNum, 2).Value = $val
    } else {
        $ws.Cells.Item($rowNum, 2).Value = [int]$val
    }
}

# Remove the autoFilter from the sheet
$ws.AutoFilterMode = $false

# Update the defined name range to match the new 2-column table (A:B instead of A:D)
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Bedingungszuweisung!`$A`$1:`$B`$107"
    }
}

# Restore view state as closely as possible
$ws.Activate()
$ws.Range("A115").Select()
$excel.ActiveWindow.ScrollRow = 108
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 400

Write-Output "done: $rowNum rows written"
